$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2:AF2").ClearContents()
$ws.Range("AH2:AI2").ClearContents()
$ws.Range("AG2").Value = 0
$ws.Range("AJ2").Value = 22006758

# Row 3
$ws.Range("D3:AF3").ClearContents()
$ws.Range("AI3").ClearContents()
$ws.Range("AG3").Value = 400
$ws.Range("AH3").Value = 1
$ws.Range("AJ3").Value = 25906758

# Row 4
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("Y4:Z4").ClearContents()
$ws.Range("D4").Value = 7476
$ws.Range("E4").Value = 584
$ws.Range("F4").Value = 584
$ws.Range("G4").Value = 688
$ws.Range("H4").Value = 530
$ws.Range("I4").Value = 530
$ws.Range("K4").Value = 6003
$ws.Range("L4").Value = 3284
$ws.Range("M4").Value = 2719
$ws.Range("N4").Value = 2719
$ws.Range("P4").Value = 1315
$ws.Range("Q4").Value = 1102
$ws.Range("R4").Value = -618
$ws.Range("S4").Value = -73
$ws.Range("T4").Value = 304
$ws.Range("U4").Value = 798
$ws.Range("V4").Value = 203
$ws.Range("W4").Value = 7.81
$ws.Range("X4").Value = 7.09
$ws.Range("AA4").Value = 120.76
$ws.Range("AB4").Value = 106.75
$ws.Range("AC4").Value = 2038
$ws.Range("AD4").Value = 12.31
$ws.Range("AE4").Value = 10342
$ws.Range("AF4").Value = 2.43
$ws.Range("AG4").Value = 500
$ws.Range("AH4").Value = 1.99
$ws.Range("AI4").Value = 24.82
$ws.Range("AJ4").Value = 26294258

# Row 5
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()
$ws.Range("D5").Value = 9964
$ws.Range("E5").Value = 1013
$ws.Range("F5").Value = 1013
$ws.Range("G5").Value = 1012
$ws.Range("H5").Value = 778
$ws.Range("I5").Value = 778
$ws.Range("K5").Value = 8004
$ws.Range("L5").Value = 4690
$ws.Range("M5").Value = 3314
$ws.Range("N5").Value = 3314
$ws.Range("P5").Value = 1318
$ws.Range("Q5").Value = 1446
$ws.Range("R5").Value = -2234
$ws.Range("S5").Value = 202
$ws.Range("T5").Value = 1332
$ws.Range("U5").Value = 114
$ws.Range("V5").Value = 545
$ws.Range("W5").Value = 10.17
$ws.Range("X5").Value = 7.81
$ws.Range("Y5").Value = 25.79
$ws.Range("Z5").Value = 11.11
$ws.Range("AA5").Value = 141.51
$ws.Range("AB5").Value = 155.2
$ws.Range("AC5").Value = 2954
$ws.Range("AD5").Value = 11.97
$ws.Range("AE5").Value = 12643
$ws.Range("AF5").Value = 2.8
$ws.Range("AG5").Value = 600
$ws.Range("AH5").Value = 1.7
$ws.Range("AI5").Value = 20.22
$ws.Range("AJ5").Value = 26356758

# Row 6
$ws.Range("AG6:AH6").ClearContents()
$ws.Range("D6").Value = 12594
$ws.Range("E6").Value = 1012
$ws.Range("F6").Value = 1012
$ws.Range("G6").Value = 933
$ws.Range("H6").Value = 709
$ws.Range("I6").Value = 709
$ws.Range("K6").Value = 10316
$ws.Range("L6").Value = 6492
$ws.Range("M6").Value = 3824
$ws.Range("N6").Value = 3824
$ws.Range("P6").Value = 1318
$ws.Range("Q6").Value = 1215
$ws.Range("R6").Value = -63
$ws.Range("S6").Value = 47
$ws.Range("T6").Value = 1535
$ws.Range("U6").Value = -320
$ws.Range("V6").Value = 1141
$ws.Range("W6").Value = 8.04
$ws.Range("X6").Value = 5.63
$ws.Range("Y6").Value = 19.86
$ws.Range("Z6").Value = 7.74
$ws.Range("AA6").Value = 169.76
$ws.Range("AB6").Value = 191.6
$ws.Range("AC6").Value = 2689
$ws.Range("AD6").Value = 12.48
$ws.Range("AE6").Value = 14540
$ws.Range("AF6").Value = 2.31
$ws.Range("AI6").Value = 24.12
$ws.Range("AJ6").Value = 26356758

# Row 7
$ws.Range("D7").Value = 13892
$ws.Range("E7").Value = -262
$ws.Range("G7").Value = -608
$ws.Range("H7").Value = -426
$ws.Range("I7").Value = -440
$ws.Range("K7").Value = 14435
$ws.Range("L7").Value = 11187
$ws.Range("M7").Value = 3247
$ws.Range("N7").Value = 3257
$ws.Range("P7").Value = 1319
$ws.Range("Q7").Value = 1666
$ws.Range("R7").Value = -3207
$ws.Range("S7").Value = 1033
$ws.Range("T7").Value = 1073
$ws.Range("U7").Value = -476
$ws.Range("W7").Value = -1.88
$ws.Range("X7").Value = -3.07
$ws.Range("Y7").Value = -12.44
$ws.Range("Z7").Value = -3.45
$ws.Range("AA7").Value = 344.5
$ws.Range("AC7").Value = -1671
$ws.Range("AD7").Value = -13.14
$ws.Range("AE7").Value = 12389
$ws.Range("AF7").Value = 1.77
$ws.Range("AG7").Value = 315
$ws.Range("AH7").Value = 1.44
$ws.Range("AI7").Value = -18.87

# Row 8
$ws.Range("D8").Value = 15077
$ws.Range("E8").Value = 213
$ws.Range("G8").Value = 89
$ws.Range("H8").Value = 66
$ws.Range("I8").Value = 60
$ws.Range("K8").Value = 14857
$ws.Range("L8").Value = 11620
$ws.Range("M8").Value = 3238
$ws.Range("N8").Value = 3320
$ws.Range("P8").Value = 1321
$ws.Range("Q8").Value = 1635
$ws.Range("R8").Value = -1185
$ws.Range("S8").Value = -67
$ws.Range("T8").Value = 918
$ws.Range("U8").Value = 514
$ws.Range("W8").Value = 1.42
$ws.Range("X8").Value = 0.44
$ws.Range("Y8").Value = 1.82
$ws.Range("Z8").Value = 0.45
$ws.Range("AA8").Value = 358.87
$ws.Range("AC8").Value = 227
$ws.Range("AD8").Value = 96.88
$ws.Range("AE8").Value = 12631
$ws.Range("AF8").Value = 1.74
$ws.Range("AG8").Value = 342
$ws.Range("AH8").Value = 1.56
$ws.Range("AI8").Value = 151.09

# Row 9
$ws.Range("D9").Value = 16520
$ws.Range("E9").Value = 744
$ws.Range("G9").Value = 633
$ws.Range("H9").Value = 474
$ws.Range("I9").Value = 473
$ws.Range("K9").Value = 15746
$ws.Range("L9").Value = 12165
$ws.Range("M9").Value = 3582
$ws.Range("N9").Value = 3763
$ws.Range("P9").Value = 1321
$ws.Range("Q9").Value = 1847
$ws.Range("R9").Value = -1447
$ws.Range("S9").Value = 78
$ws.Range("T9").Value = 1178
$ws.Range("U9").Value = 511
$ws.Range("W9").Value = 4.51
$ws.Range("X9").Value = 2.87
$ws.Range("Y9").Value = 13.36
$ws.Range("Z9").Value = 3.1
$ws.Range("AA9").Value = 339.64
$ws.Range("AC9").Value = 1795
$ws.Range("AD9").Value = 12.23
$ws.Range("AE9").Value = 14314
$ws.Range("AF9").Value = 1.53
$ws.Range("AG9").Value = 488
$ws.Range("AH9").Value = 2.23
$ws.Range("AI9").Value = 27.21
